$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update cell values (column B and C) ---
$ws.Range("B2").Value = 15
$ws.Range("C2").Value = 15
$ws.Range("B3").Value = 5776764.6
$ws.Range("B4").Value = 47266.3
$ws.Range("B5").Value = 8685.6379
$ws.Range("B6").Value = 5.7542947
$ws.Range("B7").Value = 10.3535155
$ws.Range("B8").Value = 1.7854613
$ws.Range("B9").Value = 34.400196
$ws.Range("B10").Value = 16.396242
$ws.Range("B11").Value = 10.4402165
$ws.Range("B12").Value = 381.348126
$ws.Range("B13").Value = 35.816368
$ws.Range("B14").Value = 458614.7000000001
$ws.Range("B15").Value = 23994.5295
$ws.Range("B16").Value = 460.37185
$ws.Range("B17").Value = 861.26382
$ws.Range("B18").Value = 45.100607
$ws.Range("B19").Value = 16.56282
$ws.Range("B20").Value = 21814.2201
$ws.Range("B21").Value = 593.92998
$ws.Range("B22").Value = 18.967324
$ws.Range("B23").Value = 118.58722
$ws.Range("B24").Value = 69.65235199999999
$ws.Range("B25").Value = 320.15789
$ws.Range("B26").Value = 156.57514
$ws.Range("B27").Value = 352.12046
$ws.Range("B28").Value = 202.888319
$ws.Range("B29").Value = 392.98692
$ws.Range("B30").Value = 331.106875
$ws.Range("B31").Value = 276.35072
$ws.Range("B32").Value = 28409.4398

# --- Swap the cell formatting (style) between column B and C for rows 9, 19, 30 ---
# Use an unused helper cell to hold a format temporarily while swapping so that
# no new style definitions are introduced (mirrors the existing style reuse).
function Swap-Format($rowNum) {
    $bCell = $ws.Range("B$rowNum")
    $cCell = $ws.Range("C$rowNum")
    $helper = $ws.Range("Z1")

    $bCell.Copy()
    $helper.PasteSpecial(-4122)  # xlPasteFormats

    $cCell.Copy()
    $bCell.PasteSpecial(-4122)

    $helper.Copy()
    $cCell.PasteSpecial(-4122)

    $helper.Clear()
    $excel.CutCopyMode = $false
}

Swap-Format 9
Swap-Format 19
Swap-Format 30
